# Regenerate orders with updated distance/size codes.
# The experiment's Distance and Size condition codes changed:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# These codes appear throughout the sheet (Condition, Filename_Left,
# Filename_Right, Distance, Size columns), so do a global find/replace
# across every cell on the active sheet for each code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Worksheet.Cells.Replace(What, Replacement, LookAt, SearchOrder, MatchCase, MatchByte, SearchFormat, ReplaceFormat)
# LookAt: 2 = xlPart (substring match, like the plain string substitutions in the source diff)
# SearchOrder: 1 = xlByRows
$ws.Cells.Replace("D64", "D69", 2, 1, $false, $false, $false, $false) | Out-Null
$ws.Cells.Replace("D51", "D55", 2, 1, $false, $false, $false, $false) | Out-Null
$ws.Cells.Replace("D80", "D86", 2, 1, $false, $false, $false, $false) | Out-Null
$ws.Cells.Replace("S30", "S31", 2, 1, $false, $false, $false, $false) | Out-Null
